$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 101
$ws.Range("D2").Value = 87
$ws.Range("E2").Value = 0.8613861386138614
$ws.Range("F2").Value = 0.8529411764705882
$ws.Range("G2").Value = 0.09777407078138732
$ws.Range("H2").Value = 0.08339553096059507
$ws.Range("I2").Value = 459306.928978237
$ws.Range("J2").Value = 167039.4646901185
$ws.Range("L2").Value = 167039.4646901185
$ws.Range("M2").Value = 626346.3936683555
$ws.Range("N2").Value = 10084600.3388
$ws.Range("O2").Value = 9676859.4087
$ws.Range("P2").Value = 0.01656381602426448
$ws.Range("Q2").Value = 0.01726174346812782

$ws.Range("C3").Value = 102
$ws.Range("D3").Value = 88
$ws.Range("E3").Value = 0.8627450980392157
$ws.Range("F3").Value = 0.8543689320388349
$ws.Range("G3").Value = 0.0981516967687382
$ws.Range("H3").Value = 0.08385776034610644
$ws.Range("I3").Value = 484236.7288196762
$ws.Range("J3").Value = 176905.9446168681
$ws.Range("L3").Value = 176905.9446168681
$ws.Range("M3").Value = 661142.6734365443
$ws.Range("N3").Value = 10481710.504064
$ws.Range("O3").Value = 10074337.346061
$ws.Range("P3").Value = 0.01687758353450781
$ws.Range("Q3").Value = 0.01756005765342344

$ws.Range("D4").Value = 88
$ws.Range("E4").Value = 0.8461538461538461
$ws.Range("F4").Value = 0.8461538461538461
$ws.Range("G4").Value = 0.09881066311551624
$ws.Range("H4").Value = 0.08360902263620605
$ws.Range("I4").Value = 510062.6018105842
$ws.Range("J4").Value = 182732.508518533
$ws.Range("L4").Value = 182732.508518533
$ws.Range("M4").Value = 692795.1103291172
$ws.Range("N4").Value = 10837517.91078592
$ws.Range("O4").Value = 10429173.55804283
$ws.Range("P4").Value = 0.01686110325471023
$ws.Range("Q4").Value = 0.0175212836857636

$ws.Range("D5").Value = 89
$ws.Range("E5").Value = 0.8476190476190476
$ws.Range("F5").Value = 0.8476190476190476
$ws.Range("G5").Value = 0.09810748281467674
$ws.Range("H5").Value = 0.08315777114767837
$ws.Range("I5").Value = 528942.5728075609
$ws.Range("J5").Value = 189598.5302454186
$ws.Range("L5").Value = 189598.5302454186
$ws.Range("M5").Value = 718541.1030529796
$ws.Range("N5").Value = 11253820.6075095
$ws.Range("O5").Value = 10843125.92418412
$ws.Range("P5").Value = 0.01684748112289106
$ws.Range("Q5").Value = 0.01748559701059497

$ws.Range("D6").Value = 89
$ws.Range("E6").Value = 0.839622641509434
$ws.Range("F6").Value = 0.839622641509434
$ws.Range("G6").Value = 0.09785518137528718
$ws.Range("H6").Value = 0.08216142587170339
$ws.Range("I6").Value = 546463.8167449427
$ws.Range("J6").Value = 195297.9695293586
$ws.Range("L6").Value = 195297.9695293586
$ws.Range("M6").Value = 741761.7862743011
$ws.Range("N6").Value = 11744042.02153478
$ws.Range("O6").Value = 11329576.49770964
$ws.Range("P6").Value = 0.01662953599546435
$ws.Range("Q6").Value = 0.01723788789182363
